# Applies the "stateless entities outside the US" edit:
# Inserts two new metric columns ("IMF (20%) - Sales" and "IMF (20%) - Sales + Emp")
# in place of the old "IMF - Sales" / "IMF - Sales + Emp" columns (F & G),
# pushes the old IMF values into the old "OECD (20%)" columns (H & I),
# and removes the "OECD (20%)" labels/columns entirely (they are replaced
# by the shifted-in IMF data). The "OECD - Sales" / "OECD - Sales + Emp"
# columns (J & K) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1) ---
$ws.Range("F1").Value = "IMF (20%) - Sales"
$ws.Range("G1").Value = "IMF (20%) - Sales + Emp"
$ws.Range("H1").Value = "IMF - Sales"
$ws.Range("I1").Value = "IMF - Sales + Emp"
# J1/K1 ("OECD - Sales" / "OECD - Sales + Emp") remain unchanged.

# --- Row 2: High Income ---
# Capture the old IMF values (currently in F2/G2) before overwriting them.
# Use .Value2 (not .Value) so the numeric value is resolved immediately.
$oldF2 = $ws.Range("F2").Value2
$oldG2 = $ws.Range("G2").Value2
$ws.Range("F2").Value = 0.007415338153038805
$ws.Range("G2").Value = 0.006034852546579005
$ws.Range("H2").Value = $oldF2
$ws.Range("I2").Value = $oldG2

# --- Row 3: LICs (all 'inf' strings) --- no numeric change needed; values already equal.

# --- Row 4: LMICs ---
$oldF4 = $ws.Range("F4").Value2
$oldG4 = $ws.Range("G4").Value2
$ws.Range("F4").Value = -0.002969016773552038
$ws.Range("G4").Value = 0.01497156856423381
$ws.Range("H4").Value = $oldF4
$ws.Range("I4").Value = $oldG4

# --- Row 5: Tax haven ---
$oldF5 = $ws.Range("F5").Value2
$oldG5 = $ws.Range("G5").Value2
$ws.Range("F5").Value = -0.1074605673108959
$ws.Range("G5").Value = -0.1308112174461001
$ws.Range("H5").Value = $oldF5
$ws.Range("I5").Value = $oldG5

# --- Row 6: UMICs ---
$oldF6 = $ws.Range("F6").Value2
$oldG6 = $ws.Range("G6").Value2
$ws.Range("F6").Value = 0.002028407217499605
$ws.Range("G6").Value = 0.00515745871754894
$ws.Range("H6").Value = $oldF6
$ws.Range("I6").Value = $oldG6

$wb.Save()
